$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New section header label "distance" (new shared string) in A27
$ws.Range("A27").Value = "distance"

# New data rows 28-32
$ws.Range("A28").Value = 0
$ws.Range("B28").Value = 30

$ws.Range("A29").Value = 0.05
$ws.Range("B29").Value = 23

$ws.Range("A30").Value = 0.1
$ws.Range("B30").Value = 20

$ws.Range("A31").Value = 0.3
$ws.Range("B31").Value = 0.5

$ws.Range("A32").Value = 0.535
$ws.Range("B32").Value = 0

# Update the selection to match the final state of the workbook
$ws.Range("A28:B32").Select()
